$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage original rows 2-13 into scratch rows 102-113 (preserve originals before overwriting,
# since the row permutation below has cycles: [3,6,7] and [8,10,11,12,13]).
for ($i = 2; $i -le 13; $i++) {
    $scratchRow = $i + 100
    $ws.Rows($i).Copy()
    $ws.Rows($scratchRow).PasteSpecial()
}

# Copy each destination row from its staged source row per the commits row permutation.
$ws.Rows(102).Copy()
$ws.Rows(2).PasteSpecial()
$ws.Rows(106).Copy()
$ws.Rows(3).PasteSpecial()
$ws.Rows(104).Copy()
$ws.Rows(4).PasteSpecial()
$ws.Rows(105).Copy()
$ws.Rows(5).PasteSpecial()
$ws.Rows(107).Copy()
$ws.Rows(6).PasteSpecial()
$ws.Rows(103).Copy()
$ws.Rows(7).PasteSpecial()
$ws.Rows(110).Copy()
$ws.Rows(8).PasteSpecial()
$ws.Rows(109).Copy()
$ws.Rows(9).PasteSpecial()
$ws.Rows(111).Copy()
$ws.Rows(10).PasteSpecial()
$ws.Rows(112).Copy()
$ws.Rows(11).PasteSpecial()
$ws.Rows(113).Copy()
$ws.Rows(12).PasteSpecial()
$ws.Rows(108).Copy()
$ws.Rows(13).PasteSpecial()

# Clear the scratch rows used for staging.
$ws.Range("A102:AY113").ClearContents()

# Update Taxonsorteringsordning (column B) values per the commit.
$ws.Range("B2").Value = 89820
$ws.Range("B3").Value = 89535
$ws.Range("B4").Value = 89539
$ws.Range("B5").Value = 77636
$ws.Range("B6").Value = 77636
$ws.Range("B7").Value = 89557
$ws.Range("B8").Value = 89535
$ws.Range("B9").Value = 89535
$ws.Range("B10").Value = 77636
$ws.Range("B11").Value = 90821
$ws.Range("B12").Value = 90221
$ws.Range("B13").Value = 77636

# Clear cells that no longer apply after the row content swap (substrate fields).
$ws.Range("AJ6").ClearContents()
$ws.Range("AK6").ClearContents()
$ws.Range("AJ11").ClearContents()
$ws.Range("AK11").ClearContents()
$ws.Range("AM11").ClearContents()
$ws.Range("AO11").ClearContents()
$ws.Range("AJ13").ClearContents()
$ws.Range("AK13").ClearContents()
$ws.Range("AM13").ClearContents()
$ws.Range("AO13").ClearContents()
